$wb = $excel.ActiveWorkbook

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35705.445
$ws.Range("I2").Value = 39668.625
$ws.Range("K2").Value = 39668.625
$ws.Range("M2").Value = -39555.625
$ws.Range("H32").Value = 16671106
$ws.Range("I32").Value = 17245906
$ws.Range("K32").Value = 17245906
$ws.Range("M32").Value = -17245619
$ws.Range("H45").Value = 3120.6667
$ws.Range("I45").Value = 2681.25
$ws.Range("K45").Value = 2681.25
$ws.Range("M45").Value = -2304.25
$ws.Range("H102").Value = 10678.75
$ws.Range("I102").Value = 10631.286
$ws.Range("K102").Value = 10631.286
$ws.Range("M102").Value = -9009.286
$ws.Range("H104").Value = 43681.25
$ws.Range("J104").Value = 43681.25
$ws.Range("L104").Value = 43681.25
$ws.Range("N104").Value = -50669.25
$ws.Range("H110").Value = 5745.7
$ws.Range("I110").Value = 4932.125
$ws.Range("J110").Value = 9000
$ws.Range("K110").Value = 4932.125
$ws.Range("L110").Value = 9000
$ws.Range("M110").Value = -2887.125
$ws.Range("N110").Value = -13090
$ws.Range("H116").Value = 35705.445
$ws.Range("I116").Value = 39668.625
$ws.Range("K116").Value = 39668.625
$ws.Range("M116").Value = -37374.625
$ws.Range("H122").Value = 2673.2
$ws.Range("I122").Value = 2414.6667
$ws.Range("K122").Value = 7244.000100000001
$ws.Range("M122").Value = -4794.000100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35705.445
$ws.Range("I3").Value = 39668.625
$ws.Range("K3").Value = 39668.625
$ws.Range("M3").Value = -39554.625
$ws.Range("H46").Value = 4000
$ws.Range("J46").Value = 4000
$ws.Range("L46").Value = 4000
$ws.Range("N46").Value = -4596
$ws.Range("H86").Value = 2991.111
$ws.Range("I86").Value = 2780.5454
$ws.Range("J86").Value = 3322
$ws.Range("K86").Value = 2780.5454
$ws.Range("L86").Value = 3322
$ws.Range("M86").Value = -1657.5454
$ws.Range("N86").Value = -5568
$ws.Range("H89").Value = 2991.111
$ws.Range("I89").Value = 2780.5454
$ws.Range("J89").Value = 3322
$ws.Range("K89").Value = 13902.727
$ws.Range("L89").Value = 16610
$ws.Range("M89").Value = -8286.726999999999
$ws.Range("N89").Value = -27842
$ws.Range("H105").Value = 1852.4375
$ws.Range("I105").Value = 1881.4286
$ws.Range("J105").Value = 1649.5
$ws.Range("K105").Value = 1881.4286
$ws.Range("L105").Value = 1649.5
$ws.Range("M105").Value = -134.4286
$ws.Range("N105").Value = -5143.5
$ws.Range("H134").Value = 1148.4359
$ws.Range("I134").Value = 952.4722
$ws.Range("K134").Value = 2857.4166
$ws.Range("M134").Value = -322.4166

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 110250
$ws.Range("I4").Value = 110250
$ws.Range("K4").Value = 110250
$ws.Range("M4").Value = -110138
$ws.Range("H58").Value = 2613.9333
$ws.Range("I58").Value = 1950.9
$ws.Range("J58").Value = 3940
$ws.Range("K58").Value = 1950.9
$ws.Range("L58").Value = 3940
$ws.Range("M58").Value = -1747.9
$ws.Range("N58").Value = -4346
$ws.Range("H132").Value = 3046.5715
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 3020.3513
$ws.Range("I134").Value = 2252.1428
$ws.Range("K134").Value = 6756.428400000001
$ws.Range("M134").Value = -4221.428400000001
$ws.Range("H136").Value = 2613.9333
$ws.Range("I136").Value = 1950.9
$ws.Range("J136").Value = 3940
$ws.Range("K136").Value = 5852.700000000001
$ws.Range("L136").Value = 11820
$ws.Range("M136").Value = -3302.700000000001
$ws.Range("N136").Value = -16920

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35893170
$ws.Range("I4").Value = 40200228
$ws.Range("K4").Value = 120600684
$ws.Range("M4").Value = -120600572
$ws.Range("H37").Value = 387134.75
$ws.Range("J37").Value = 387134.75
$ws.Range("L37").Value = 1161404.25
$ws.Range("N37").Value = -1161628.25
$ws.Range("H80").Value = 3888.7778
$ws.Range("I80").Value = 3333.1667
$ws.Range("K80").Value = 9999.500100000001
$ws.Range("M80").Value = -9063.500100000001
$ws.Range("H83").Value = 3888.7778
$ws.Range("I83").Value = 3333.1667
$ws.Range("K83").Value = 29998.5003
$ws.Range("M83").Value = -25318.5003
$ws.Range("H97").Value = 400
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1200
$ws.Range("N97").Value = -2192
$ws.Range("M97").ClearContents()
$ws.Range("H132").Value = 2044.3334
$ws.Range("I132").Value = 1857
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 16713
$ws.Range("L132").Value = 24300
$ws.Range("M132").Value = -14183
$ws.Range("N132").Value = -29360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 100261
$ws.Range("J39").Value = 100261
$ws.Range("L39").Value = 100261
$ws.Range("N39").Value = -101325
$ws.Range("H70").Value = 188319.33
$ws.Range("I70").Value = 280744.75
$ws.Range("J70").Value = 3468.5
$ws.Range("K70").Value = 280744.75
$ws.Range("L70").Value = 3468.5
$ws.Range("M70").Value = -280474.75
$ws.Range("N70").Value = -4008.5
$ws.Range("H73").Value = 188319.33
$ws.Range("I73").Value = 280744.75
$ws.Range("J73").Value = 3468.5
$ws.Range("K73").Value = 280744.75
$ws.Range("L73").Value = 3468.5
$ws.Range("M73").Value = -279808.75
$ws.Range("N73").Value = -5340.5
$ws.Range("H102").Value = 3322.5667
$ws.Range("I102").Value = 2258.158
$ws.Range("K102").Value = 2258.158
$ws.Range("M102").Value = -636.1579999999999
$ws.Range("H122").Value = 10285.429
$ws.Range("I122").Value = 17000
$ws.Range("K122").Value = 51000
$ws.Range("M122").Value = -48550
$ws.Range("H132").Value = 3197.2222
$ws.Range("I132").Value = 3197.2222
$ws.Range("K132").Value = 9591.6666
$ws.Range("M132").Value = -7061.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 14812.5
$ws.Range("J20").Value = 14785.714
$ws.Range("L20").Value = 14785.714
$ws.Range("N20").Value = -15237.714
$ws.Range("H122").Value = 7653.846
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7653.846
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 22961.538
$ws.Range("N122").Value = -27861.538
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 6517.3335
$ws.Range("I132").Value = 2820.8
$ws.Range("J132").Value = 25000
$ws.Range("K132").Value = 8462.400000000001
$ws.Range("L132").Value = 75000
$ws.Range("M132").Value = -5932.400000000001
$ws.Range("N132").Value = -80060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2749.3333
$ws.Range("I132").Value = 3166.3333
$ws.Range("J132").Value = 2332.3333
$ws.Range("K132").Value = 9498.999899999999
$ws.Range("L132").Value = 6996.999899999999
$ws.Range("M132").Value = -6968.999899999999
$ws.Range("N132").Value = -12056.9999
